# Week-04 deck: bump the cached "datetimeFigureOut" footer-date field
# (auto-populated by PowerPoint's Header & Footer dialog / whole-deck
# resave) from 8/2/2022 to 8/3/2022 everywhere it is cached: the slide
# master, every slide layout, and the notes master.

$p = $ppt.ActivePresentation

$oldDate = "8/2/2022"
$newDate = "8/3/2022"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes
